$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 745.4545000000001
$ws.Range("I19").Value = 630
$ws.Range("J19").Value = 841.6667
$ws.Range("K19").Value = 630
$ws.Range("L19").Value = 841.6667
$ws.Range("M19").Value = -455
$ws.Range("N19").Value = -1191.6667

$ws.Range("H51").Value = 2021.625
$ws.Range("I51").Value = 1133.3334
$ws.Range("J51").Value = 2554.6
$ws.Range("K51").Value = 1133.3334
$ws.Range("L51").Value = 2554.6
$ws.Range("M51").Value = -649.3334
$ws.Range("N51").Value = -3522.6

$ws.Range("H138").Value = 5413741
$ws.Range("I138").Value = 2356102.8
$ws.Range("J138").Value = 6413354
$ws.Range("K138").Value = 7068308.399999999
$ws.Range("L138").Value = 19240062
$ws.Range("M138").Value = -7063168.399999999
$ws.Range("N138").Value = -19250342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2689.0693
$ws.Range("I32").Value = 2788.9048
$ws.Range("J32").Value = 1990.2222
$ws.Range("K32").Value = 2788.9048
$ws.Range("L32").Value = 1990.2222
$ws.Range("M32").Value = -2501.9048
$ws.Range("N32").Value = -2564.2222

$ws.Range("H74").Value = 8409.315000000001
$ws.Range("I74").Value = 2579.2222
$ws.Range("J74").Value = 13656.4
$ws.Range("K74").Value = 2579.2222
$ws.Range("L74").Value = 13656.4
$ws.Range("M74").Value = -1705.2222
$ws.Range("N74").Value = -15404.4

$ws.Range("H77").Value = 8409.315000000001
$ws.Range("I77").Value = 2579.2222
$ws.Range("J77").Value = 13656.4
$ws.Range("K77").Value = 12896.111
$ws.Range("L77").Value = 68282
$ws.Range("M77").Value = -8528.111000000001
$ws.Range("N77").Value = -77018

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4392.364
$ws.Range("I86").Value = 1336.5652
$ws.Range("J86").Value = 11420.7
$ws.Range("K86").Value = 1336.5652
$ws.Range("L86").Value = 11420.7
$ws.Range("M86").Value = -213.5652
$ws.Range("N86").Value = -13666.7

$ws.Range("H89").Value = 4392.364
$ws.Range("I89").Value = 1336.5652
$ws.Range("J89").Value = 11420.7
$ws.Range("K89").Value = 6682.826
$ws.Range("L89").Value = 57103.5
$ws.Range("M89").Value = -1066.826
$ws.Range("N89").Value = -68335.5

$ws.Range("H105").Value = 3236.1428
$ws.Range("I105").Value = 3189.611
$ws.Range("K105").Value = 3189.611
$ws.Range("M105").Value = -1442.611

$ws.Range("H107").Value = 841.6957
$ws.Range("I107").Value = 885.06665
$ws.Range("J107").Value = 760.375
$ws.Range("K107").Value = 885.06665
$ws.Range("L107").Value = 760.375
$ws.Range("M107").Value = 1034.93335
$ws.Range("N107").Value = -4600.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1419.6786
$ws.Range("I31").Value = 980.1111
$ws.Range("K31").Value = 980.1111
$ws.Range("M31").Value = -685.1111

$ws.Range("H34").Value = 1419.6786
$ws.Range("I34").Value = 980.1111
$ws.Range("K34").Value = 980.1111
$ws.Range("M34").Value = -778.1111

$ws.Range("H58").Value = 1688.6923
$ws.Range("I58").Value = 1028.3636
$ws.Range("J58").Value = 2172.9333
$ws.Range("K58").Value = 1028.3636
$ws.Range("L58").Value = 2172.9333
$ws.Range("M58").Value = -825.3635999999999
$ws.Range("N58").Value = -2578.9333

$ws.Range("H132").Value = 1916.7413
$ws.Range("I132").Value = 1528.0212
$ws.Range("K132").Value = 4584.063599999999
$ws.Range("M132").Value = -2054.063599999999

$ws.Range("H136").Value = 1688.6923
$ws.Range("I136").Value = 1028.3636
$ws.Range("J136").Value = 2172.9333
$ws.Range("K136").Value = 3085.0908
$ws.Range("L136").Value = 6518.7999
$ws.Range("M136").Value = -535.0907999999999
$ws.Range("N136").Value = -11618.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 15152787
$ws.Range("I132").Value = 797.75
$ws.Range("J132").Value = 23811066
$ws.Range("K132").Value = 7179.75
$ws.Range("L132").Value = 214299594
$ws.Range("M132").Value = -4649.75
$ws.Range("N132").Value = -214304654

$ws.Range("H137").Value = 6316022
$ws.Range("I137").Value = 8336110
$ws.Range("J137").Value = 255758.25
$ws.Range("K137").Value = 25008330
$ws.Range("L137").Value = 767274.75
$ws.Range("M137").Value = -25003230
$ws.Range("N137").Value = -777474.75

$ws.Range("H141").Value = 4058.0908
$ws.Range("I141").Value = 4642.375
$ws.Range("K141").Value = 13927.125
$ws.Range("M141").Value = -8747.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5825.4243
$ws.Range("I70").Value = 6225.1904
$ws.Range("J70").Value = 5125.8335
$ws.Range("K70").Value = 6225.1904
$ws.Range("L70").Value = 5125.8335
$ws.Range("M70").Value = -5955.1904
$ws.Range("N70").Value = -5665.8335

$ws.Range("H73").Value = 5825.4243
$ws.Range("I73").Value = 6225.1904
$ws.Range("J73").Value = 5125.8335
$ws.Range("K73").Value = 6225.1904
$ws.Range("L73").Value = 5125.8335
$ws.Range("M73").Value = -5289.1904
$ws.Range("N73").Value = -6997.8335

$ws.Range("H98").Value = 100643
$ws.Range("J98").Value = 100643
$ws.Range("L98").Value = 100643
$ws.Range("N98").Value = -106633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3177.8696
$ws.Range("I7").Value = 2000.3334
$ws.Range("J7").Value = 3354.5
$ws.Range("K7").Value = 2000.3334
$ws.Range("L7").Value = 3354.5
$ws.Range("M7").Value = -1888.3334
$ws.Range("N7").Value = -3578.5

$ws.Range("H22").Value = 756.2381
$ws.Range("I22").Value = 473.75
$ws.Range("K22").Value = 473.75
$ws.Range("M22").Value = -178.75

$ws.Range("H27").Value = 756.2381
$ws.Range("I27").Value = 473.75
$ws.Range("K27").Value = 473.75
$ws.Range("M27").Value = -366.75

$ws.Range("H40").Value = 4195
$ws.Range("I40").Value = 2466.6667
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 2466.6667
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -2330.6667
$ws.Range("N40").Value = -4772

$ws.Range("H68").Value = 2143.0588
$ws.Range("I68").Value = 2027
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2027
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1278
$ws.Range("N68").Value = -5498

$ws.Range("H71").Value = 2143.0588
$ws.Range("I71").Value = 2027
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 10135
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -6391
$ws.Range("N71").Value = -27488

$ws.Range("H82").Value = 1062.3334
$ws.Range("I82").Value = 885.7143
$ws.Range("J82").Value = 1216.875
$ws.Range("K82").Value = 885.7143
$ws.Range("L82").Value = 1216.875
$ws.Range("M82").Value = -524.7143
$ws.Range("N82").Value = -1938.875

$ws.Range("H85").Value = 1062.3334
$ws.Range("I85").Value = 885.7143
$ws.Range("J85").Value = 1216.875
$ws.Range("K85").Value = 885.7143
$ws.Range("L85").Value = 1216.875
$ws.Range("M85").Value = 362.2857
$ws.Range("N85").Value = -3712.875

$ws.Range("H126").Value = 3177.8696
$ws.Range("I126").Value = 2000.3334
$ws.Range("J126").Value = 3354.5
$ws.Range("K126").Value = 6001.0002
$ws.Range("L126").Value = 10063.5
$ws.Range("M126").Value = -3531.0002
$ws.Range("N126").Value = -15003.5
